$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2020" column (R), one cell per data row, mirroring the value
# and the cell format of the existing "2019" column (Q) in the same row.
$newCol = [ordered]@{
    4  = 2020
    5  = 2.1
    6  = 2.4
    7  = 1.4
    8  = 3.2
    9  = 2.4
    10 = 0.8
    11 = 2.2000000000000002
    12 = 4.5
    13 = 1.4
    14 = 3.2
}

foreach ($row in $newCol.Keys) {
    $value = $newCol[$row]
    $srcCell = $ws.Range("Q$row")
    $dstCell = $ws.Range("R$row")

    # Copy the source cell's format (number format, font, borders, ...) onto
    # the destination cell before writing the value.
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)

    $dstCell.Value = $value
}

$excel.CutCopyMode = $false

# Move / resize the current selection to match the new state of the sheet.
$ws.Range("R16:R17").Select() | Out-Null
